# Updated cryptos list values (price + volume columns) per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.401.45"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.42"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.91"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6331"
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07594"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.534.27"
$ws.Range("E11").Value = "  +36.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07729"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.640.85"
$ws.Range("E13").Value = "  +25.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.984"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.92"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009930"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.172"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.419.02"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.54"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.610"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.55"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1395"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.455"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.70"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05821"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.263"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.870"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.160"
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7179"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("B37").Value = "RocketPoolETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.613.63"
$ws.Range("E37").Value = "  +29.70%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.249.99"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.789"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01810"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9053"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.076"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.50"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.28"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.159"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4014"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.691"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1123"
$ws.Range("E51").Value = "  -0.38%  "
